# actualizacion Vo.Bo. UPP 4T 2020 SIPOT
# Updates the reporting period (3T -> 4T 2020) date values and refreshes
# the "Vo. Bo." sign-off: selection moves back to the top of the sheet,
# and the note cell (L8) loses its forced vertical-centering so it reads
# like the rest of the justified/wrapped notes column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 data: shift the reporting period from 3er trimestre (Jul-Sep)
#     to 4to trimestre (Oct-Dic) 2020, and bump the validation/update dates
#     (10-oct-2020 -> 10-ene-2021) to match the new period. ---
$ws.Range("B8").Value = 44105   # Fecha de inicio del periodo: 2020-10-01
$ws.Range("C8").Value = 44196   # Fecha de término del periodo: 2020-12-31
$ws.Range("J8").Value = 44206   # Fecha de validación: 2021-01-10
$ws.Range("K8").Value = 44206   # Fecha de actualización: 2021-01-10

# --- L8 ("Nota") keeps justify + wrap, but is no longer vertically
#     centered within the tall (90pt) row. ---
$ws.Range("L8").VerticalAlignment = -4107   # xlBottom (i.e. "not centered")

# --- Reset the view: scroll back to the top-left and select A3:C3 (the
#     "NOMBRE CORTO" merged header cell) instead of leaving the cursor on
#     the last-edited cell L8. ---
$ws.Range("A3:C3").Select()
